# Auto-generated: update coin price/volume table to the latest scrape
# (commit: "Updated cryptos list on Sat Nov 11 04:25:31 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '37.071.41'
$ws.Range("E2").Value = '  +0.93%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '2.043.65'
$ws.Range("E3").Value = '  -3.46%  '

# Row 4 (TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.09%  '

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.60'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.05%  '

# Row 6 (XRP)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.650'
$ws.Range("D6").ClearFormats()

# Row 7 (USDC)
$ws.Range("E7").Value = '  +0.05%  '

# Row 8 (Solana)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.11'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +16.89%  '

# Row 9 (OKB)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.43'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.21%  '

# Row 10 (Cardano)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.374'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.09%  '

# Row 11 (Dogecoin)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0760'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.66%  '

# Row 12 (TRON)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.105'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.96%  '

# Row 13 (Chainlink)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.89'
$ws.Range("D13").ClearFormats()

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = '2.344.18'
$ws.Range("E14").Value = '  -3.33%  '

# Row 15 (Polygon)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.809'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.53%  '

# Row 16 (Polkadot)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.65%  '

# Row 17 (WrappedEther)
$ws.Range("D17").Value = '2.045.96'
$ws.Range("E17").Value = '  -3.29%  '

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = '36.990.55'
$ws.Range("E18").Value = '  +0.79%  '

# Row 19: 'Litecoin' -> 'ShibaInu'
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0899'
$ws.Range("E19").Value = '  +6.95%  '

# Row 20: 'ShibaInu' -> 'Litecoin'
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.66'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.87%  '

# Row 21 (Avalanche)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.14'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.39%  '

# Row 22 (BitcoinCash)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.04'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.56%  '

# Row 23 (Uniswap)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.22'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.00%  '

# Row 24 (Dai)
$ws.Range("E24").Value = '  -0.05%  '

# Row 25 (Toncoin)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.61%  '

# Row 26 (Monero)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.44%  '

# Row 27 (Cosmos)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.31%  '

# Row 28 (EthereumClassic)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.93'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.70%  '

# Row 29 (PancakeSwap)
$ws.Range("E29").Value = '  -2.79%  '

# Row 30 (Stellar)
$ws.Range("E30").Value = '  -1.26%  '

# Row 31 (Filecoin)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.53'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.26%  '

# Row 32 (ImmutableX)
$ws.Range("E32").Value = '  +12.66%  '

# Row 33 (Hedera)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0616'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.51%  '

# Row 34 (InternetComputer(DFINITY))
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.30'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.64%  '

# Row 35: 'BinanceUSD' -> 'Kaspa'
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0881'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.04%  '

# Row 36: 'Kaspa' -> 'BinanceUSD'
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.06%  '

# Row 37 (LidoDAOToken)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.24'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.91%  '

# Row 38 (WEMIXToken)
$ws.Range("E38").Value = '  -6.29%  '

# Row 39 (TrustWalletToken)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.32'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.85%  '

# Row 40: 'Gas' -> 'Cronos'
$ws.Range("B40").Value = 'Cronos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.102'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +20.69%  '

# Row 41: 'Cronos' -> 'Gas'
$ws.Range("B41").Value = 'Gas'
$ws.Range("C41").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.87'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -39.73%  '

# Row 42 (InjectiveProtocol)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +10.79%  '

# Row 43 (VeChain)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0221'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.08%  '

# Row 44 (ARBITRUM)
$ws.Range("E44").Value = '  -5.35%  '

# Row 45 (Aave)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.09'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.07%  '

# Row 46 (HuobiToken)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.33%  '

# Row 47 (FTXToken)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.04'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +43.75%  '

# Row 48 (Maker)
$ws.Range("D48").Value = '1.290.13'
$ws.Range("E48").Value = '  -5.06%  '

# Row 49 (MXToken)
$ws.Range("E49").Value = '  +2.68%  '

# Row 50 (RenderToken)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.46%  '

# Row 51 (FraxShare)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.73'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.42%  '
